$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-17
$data = @(
    @(7,7),
    @(9,9),
    @(7,8),
    @(6,7),
    @(6,7),
    @(8,8),
    @(7,8),
    @(8,8),
    @(4,6),
    @(8,8),
    @(9,9),
    @(6,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
